$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" data row (originally row 26) entirely.
$ws.Rows(26).Delete()

# After the above deletion, the row that used to be "SC 92" (originally
# row 28) has shifted up to row 27. Remove it too.
$ws.Rows(27).Delete()

# With both rows removed, the remaining rows have shifted up by either one
# or two positions. Update column E (the "D" header column) values for the
# surviving rows to their corrected values.
$ws.Range("E26").Value = -5        # SC 5
$ws.Range("E27").ClearContents()   # SC 101
$ws.Range("E28").ClearContents()   # SC 105
$ws.Range("E29").Value = -6.8      # SC 119
$ws.Range("E30").Value = -5.7      # SC 120
$ws.Range("E31").ClearContents()   # SC 132
$ws.Range("E32").ClearContents()   # SC 193
$ws.Range("E33").Value = -10.7     # SC 232
